$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guardian")

# Update the Guardian sheet header row: collapse the spaced header labels into
# single "PascalCase" tokens (MRN is unchanged).
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("D1").Value = "AccountNumber"
$ws.Range("E1").Value = "PhoneNumber"
$ws.Range("F1").Value = "CellphoneNumber"
$ws.Range("G1").Value = "EmailAddress"

# Move the active tab/selection from Patient to Guardian.
$ws.Activate() | Out-Null
$ws.Range("G1").Select() | Out-Null
